$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.105.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.853.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.51%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Value = "'237.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.16%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.75%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.21%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'41.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.32%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.27%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.124.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.49%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.848.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.13%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +1.54%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'35.092.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.62%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'69.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.50%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0793"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.91%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'240.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.36%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.55%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.14%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'167.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +23.90%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.81%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E31").Value = "'  +0.92%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.98%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +26.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.36%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +19.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +11.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +6.23%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'VeChain"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.0201"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Aave"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'90.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.36%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.338.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'14.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.92%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.48%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'12.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +48.74%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.70%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0555"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +6.32%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.52%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'6.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.74%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.036.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.39%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.23%  "
$ws.Range("E51").Style = "Normal"
